$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.002") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.819.70"
$ws.Range("E2").Value = "  -3.14%  "

$ws.Range("D3").Value = "1.798.64"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "315.64"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").Value = "0.5374"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "0.3857"
$ws.Range("E8").Value = "  +1.57%  "

$ws.Range("D9").Value = "0.07439"
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("D10").Value = "41.51"
$ws.Range("E10").Value = "  -2.82%  "

$ws.Range("D11").Value = "1.087"
$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").Value = "6.206"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").Value = "7.503"
$ws.Range("E14").Value = "  +1.54%  "

$ws.Range("D15").Value = "20.34"
$ws.Range("E15").Value = "  -3.23%  "

$ws.Range("D16").Value = "1.787.67"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").Value = "88.51"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "0.00001059"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "0.06495"
$ws.Range("E19").Value = "  +0.74%  "

$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").Value = "5.939"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "27.824.57"
$ws.Range("E23").Value = "  -3.05%  "

$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").Value = "156.15"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("D27").Value = "20.30"
$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("D28").Value = "1.997.86"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").Value = "2.332"
$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").Value = "121.60"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").Value = "1.117"
$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("D32").Value = "0.1094"
$ws.Range("E32").Value = "  +4.80%  "

$ws.Range("D33").Value = "3.651"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").Value = "5.512"
$ws.Range("E34").Value = "  -3.09%  "

$ws.Range("D35").Value = "0.07030"
$ws.Range("E35").Value = "  +8.18%  "

$ws.Range("D36").Value = "0.2200"
$ws.Range("E36").Value = "  -3.43%  "

$ws.Range("D37").Value = "0.02276"
$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("D38").Value = "5.059"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "11.33"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "8.476"
$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("D41").Value = "0.6096"
$ws.Range("E41").Value = "  -2.81%  "

$ws.Range("D42").Value = "1.414"
$ws.Range("E42").Value = "  +1.40%  "

$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("D44").Value = "13.27"
$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("D45").Value = "3.677"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").Value = "0.5700"
$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("D47").Value = "124.83"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").Value = "1.175"
$ws.Range("E48").Value = "  +1.55%  "

$ws.Range("D49").Value = "1.910"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("D50").Value = "0.06795"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").Value = "71.62"
$ws.Range("E51").Value = "  -1.74%  "

# Restore the default (Normal) style on column D so no extra number-format
# style is left attached to the cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"

